$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.633.79'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '2.615.53'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '516.07'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.71'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  +2.29%  '
$ws.Range('D9').Value = '2.627.78'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.70'
$ws.Range('E10').Value = '  +3.93%  '
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('E12').Value = '  +1.57%  '
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').Value = '3.072.57'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = '60.704.17'
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.76'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('D18').Value = '2.624.53'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.76'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '358.01'
$ws.Range('E20').Value = '  +4.54%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.67'
$ws.Range('E21').Value = '  +2.94%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.22'
$ws.Range('E22').Value = '  +1.85%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '61.17'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.427'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('D26').Value = '2.736.03'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  +1.15%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').Value = '0.0₃0846'
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.35'
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.48'
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.94'
$ws.Range('E34').Value = '  +4.38%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '150.86'
$ws.Range('E35').Value = '  -3.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.03'
$ws.Range('E36').Value = '  +1.60%  '
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.892'
$ws.Range('E38').Value = '  +6.65%  '
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.850'
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.34'
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '292.48'
$ws.Range('E43').Value = '  -2.98%  '
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.997'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0557'
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '19.72'
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.97'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0236'
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('E51').Value = '  +0.10%  '
